$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the data between row 2 and row 3 for the changed columns
# (D: Fecha, L: Calidad, N: Precio mínimo, O: Precio máximo,
#  P: Precio promedio ponderado, R: Origen, S: Precio $/Kg)

$ws.Range("D2").Value = 45043
$ws.Range("L2").Value = "Primera"
$ws.Range("N2").Value = 19000
$ws.Range("O2").Value = 20000
$ws.Range("P2").Value = 19500
$ws.Range("R2").Value = "Región de O'Higgins"
$ws.Range("S2").Value = 1083

$ws.Range("D3").Value = 45086
$ws.Range("L3").Value = "Segunda"
$ws.Range("N3").Value = 20000
$ws.Range("O3").Value = 21000
$ws.Range("P3").Value = 20500
$ws.Range("R3").Value = "Provincia de Curicó"
$ws.Range("S3").Value = 1139
